$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Balarama Holness
$ws.Range("C2").Value = 65
$ws.Range("D2").Value = 124
$ws.Range("E2").Value = 34.39153439153439
$ws.Range("F2").Value = 65.60846560846561

# Row 3 - Dexter Xurukulasuriya
$ws.Range("C3").Value = 148
$ws.Range("D3").Value = 86
$ws.Range("E3").Value = 63.24786324786324
$ws.Range("F3").Value = 36.75213675213676

# Row 5 - James Oscar
$ws.Range("C5").Value = 79
$ws.Range("D5").Value = 71
$ws.Range("E5").Value = 52.66666666666666
$ws.Range("F5").Value = 47.33333333333334

# Row 7 - Overall Average
$ws.Range("E7").Value = 47.79332615715823
$ws.Range("F7").Value = 52.20667384284177
